$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "28.100.55"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +1.57%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.791.22"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +1.77%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.002"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.37%  "

$ws.Range("E5").Value = "  -0.70%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.001"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.20%  "

$ws.Range("E7").Value = "  -4.09%  "

$ws.Range("E8").Value = "  -2.60%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "44.53"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -2.51%  "

$ws.Range("E10").Value = "  -3.41%  "

$ws.Range("E11").Value = "  -1.24%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.0000"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.28%  "

$ws.Range("E13").Value = "  -0.22%  "

$ws.Range("E14").Value = "  -0.61%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.355"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.42%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.780.85"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +1.16%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "91.85"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.95%  "

$ws.Range("E18").Value = "  -1.38%  "

$ws.Range("E19").Value = "  +1.28%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.9997"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.27%  "

$ws.Range("E21").Value = "  -1.15%  "

$ws.Range("E22").Value = "  -3.75%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "28.100.90"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +1.48%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "11.43"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -1.98%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.160"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -7.48%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "160.34"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +4.14%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "20.39"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -2.03%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.986.81"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +1.42%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.184"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -7.21%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "126.91"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -1.78%  "

$ws.Range("E31").Value = "  -3.66%  "

$ws.Range("E32").Value = "  -0.87%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.09016"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -2.79%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.503"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -5.23%  "

$ws.Range("E35").Value = "  -0.82%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "5.103"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.09%  "

$ws.Range("E38").Value = "  -0.53%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.2120"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -3.13%  "

$ws.Range("E40").Value = "  -0.98%  "

$ws.Range("E41").Value = "  -0.71%  "

$ws.Range("E42").Value = "  +0.12%  "

$ws.Range("E43").Value = "  -0.32%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "7.898"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -1.83%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "13.72"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.10%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.6000"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.22%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.710"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -1.10%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "124.48"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.14%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.994"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.45%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.154"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.76%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.06959"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.84%  "
